# DIPDAP-BOM.xlsx update
# - Replaces the obsolete X1 connector part (TSW-107-02-S-S / -RA, Alps)
#   with the new Samtech TSW-107-08-G-S-RA part, filling in the
#   previously-empty Mfg/MfgPn/Disti/Order columns for that row.
# - Updates the sheet view (zoom + selection) to match the re-saved file.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 34 ("X1" / THROUGH-HOLE connector) - new part: Samtech TSW-107-08-G-S-RA
# Set the Mfg (F) / MfgPn (G) columns first so the new shared strings are
# interned in the same order as the target workbook.
$ws.Range("F34").Value = "Samtech"
$ws.Range("G34").Value = "TSW-107-08-G-S-RA"

# Part / Device / Value columns all reference the new connector part too.
$ws.Range("B34").Value = "TSW-107-08-G-S-RA"
$ws.Range("C34").Value = "TSW-107-08-G-S-RA"
$ws.Range("D34").Value = "TSW-107-08-G-S-RA"

# Disti + order number, newly populated for this row.
$ws.Range("H34").Value = "Farnell"
$ws.Range("I34").Value = 2041451

# View state: zoomed to 150% with H42 selected.
$excel.ActiveWindow.Zoom = 150
$ws.Range("H42").Select()
